$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4995
$ws.Range("I51").Value = 4996
$ws.Range("J51").Value = 4994
$ws.Range("K51").Value = 4996
$ws.Range("L51").Value = 4994
$ws.Range("M51").Value = -4512
$ws.Range("N51").Value = -5962

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7938.5
$ws.Range("I106").Value = 7918
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 7918
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -7287
$ws.Range("N106").Value = -9262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1955.6111
$ws.Range("J121").Value = 1955.6111
$ws.Range("L121").Value = 5866.8333
$ws.Range("N121").Value = -9360.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1999.5
$ws.Range("I125").Value = 1999.5
$ws.Range("K125").Value = 17995.5
$ws.Range("M125").Value = -15535.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1014.88
$ws.Range("I132").Value = 1014.88
$ws.Range("K132").Value = 3044.64
$ws.Range("M132").Value = -514.6399999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1908.5862
$ws.Range("J137").Value = 2569
$ws.Range("L137").Value = 7707
$ws.Range("N137").Value = -12807

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3795.5386
$ws.Range("J138").Value = 3910.0303
$ws.Range("L138").Value = 11730.0909
$ws.Range("N138").Value = -22010.0909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1884.5
$ws.Range("I4").Value = 1872
$ws.Range("J4").Value = 1897
$ws.Range("K4").Value = 1872
$ws.Range("L4").Value = 1897
$ws.Range("M4").Value = -1756
$ws.Range("N4").Value = -2129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4855.048
$ws.Range("I32").Value = 4918.839
$ws.Range("K32").Value = 4918.839
$ws.Range("M32").Value = -4631.839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2102.2144
$ws.Range("I61").Value = 1323.8334
$ws.Range("K61").Value = 1323.8334
$ws.Range("M61").Value = -1111.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1377.6428
$ws.Range("I74").Value = 1053.3636
$ws.Range("J74").Value = 2566.6667
$ws.Range("K74").Value = 1053.3636
$ws.Range("L74").Value = 2566.6667
$ws.Range("M74").Value = -179.3635999999999
$ws.Range("N74").Value = -4314.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1377.6428
$ws.Range("I77").Value = 1053.3636
$ws.Range("J77").Value = 2566.6667
$ws.Range("K77").Value = 5266.817999999999
$ws.Range("L77").Value = 12833.3335
$ws.Range("M77").Value = -898.8179999999993
$ws.Range("N77").Value = -21569.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 2300.3076
$ws.Range("J132").Value = 3874.5
$ws.Range("K132").Value = 6900.9228
$ws.Range("L132").Value = 11623.5
$ws.Range("M132").Value = -4370.9228
$ws.Range("N132").Value = -16683.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2102.2144
$ws.Range("I136").Value = 1323.8334
$ws.Range("K136").Value = 3971.5002
$ws.Range("M136").Value = -1421.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2857.1428
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -753

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1943.9286
$ws.Range("I86").Value = 1708.8462
$ws.Range("K86").Value = 1708.8462
$ws.Range("M86").Value = -585.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1943.9286
$ws.Range("I89").Value = 1708.8462
$ws.Range("K89").Value = 8544.231
$ws.Range("M89").Value = -2928.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 818.2
$ws.Range("I99").Value = 774.75
$ws.Range("K99").Value = 774.75
$ws.Range("M99").Value = 723.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3047.3
$ws.Range("I134").Value = 3162.1667
$ws.Range("K134").Value = 9486.500100000001
$ws.Range("M134").Value = -6951.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1769.8334
$ws.Range("I16").Value = 1679.4
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = 1679.4
$ws.Range("L16").Value = 2222
$ws.Range("M16").Value = -1392.4
$ws.Range("N16").Value = -2796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1743.5
$ws.Range("I62").Value = 1743.5
$ws.Range("K62").Value = 1743.5
$ws.Range("M62").Value = -1119.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 1743.5
$ws.Range("I65").Value = 1743.5
$ws.Range("K65").Value = 8717.5
$ws.Range("M65").Value = -5597.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1769.8334
$ws.Range("I113").Value = 1679.4
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 1679.4
$ws.Range("L113").Value = 2222
$ws.Range("M113").Value = 490.5999999999999
$ws.Range("N113").Value = -6562

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 1700
$ws.Range("J21").Value = 1700
$ws.Range("L21").Value = 5100
$ws.Range("N21").Value = -5446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 13333.333
$ws.Range("J39").Value = 13333.333
$ws.Range("L39").Value = 39999.999
$ws.Range("N39").Value = -40587.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7237.6665
$ws.Range("J55").Value = 11890
$ws.Range("L55").Value = 35670
$ws.Range("N55").Value = -36024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108.333336
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 75000
$ws.Range("J53").Value = 75000
$ws.Range("L53").Value = 75000
$ws.Range("N53").Value = -76262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8428.857
$ws.Range("I80").Value = 7999.2
$ws.Range("J80").Value = 9503
$ws.Range("K80").Value = 7999.2
$ws.Range("L80").Value = 9503
$ws.Range("M80").Value = -7001.2
$ws.Range("N80").Value = -11499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8428.857
$ws.Range("I83").Value = 7999.2
$ws.Range("J83").Value = 9503
$ws.Range("K83").Value = 39996
$ws.Range("L83").Value = 47515
$ws.Range("M83").Value = -35004
$ws.Range("N83").Value = -57499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 392.2069
$ws.Range("I97").Value = 427.47827
$ws.Range("J97").Value = 257
$ws.Range("K97").Value = 427.47827
$ws.Range("L97").Value = 257
$ws.Range("M97").Value = 68.52172999999999
$ws.Range("N97").Value = -1249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2499.75
$ws.Range("J102").Value = 2999.5
$ws.Range("L102").Value = 2999.5
$ws.Range("N102").Value = -6243.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4167.6665
$ws.Range("I132").Value = 4004
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 12012
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -9482
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1510.8
$ws.Range("I7").Value = 1263.5
$ws.Range("K7").Value = 1263.5
$ws.Range("M7").Value = -1151.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1751
$ws.Range("I40").Value = 1751
$ws.Range("K40").Value = 1751
$ws.Range("M40").Value = -1615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1266.7142
$ws.Range("I82").Value = 1273.2
$ws.Range("J82").Value = 1250.5
$ws.Range("K82").Value = 1273.2
$ws.Range("L82").Value = 1250.5
$ws.Range("M82").Value = -912.2
$ws.Range("N82").Value = -1972.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1266.7142
$ws.Range("I85").Value = 1273.2
$ws.Range("J85").Value = 1250.5
$ws.Range("K85").Value = 1273.2
$ws.Range("L85").Value = 1250.5
$ws.Range("M85").Value = -25.20000000000005
$ws.Range("N85").Value = -3746.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1510.8
$ws.Range("I126").Value = 1263.5
$ws.Range("K126").Value = 3790.5
$ws.Range("M126").Value = -1320.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5564.2
$ws.Range("I132").Value = 5274.3335
$ws.Range("K132").Value = 15823.0005
$ws.Range("M132").Value = -13293.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3817.3076
$ws.Range("I136").Value = 4134.1
$ws.Range("K136").Value = 12402.3
$ws.Range("M136").Value = -9852.300000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9825.857
$ws.Range("I81").Value = 7262.3335
$ws.Range("K81").Value = 14524.667
$ws.Range("M81").Value = -13463.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 9825.857
$ws.Range("I84").Value = 7262.3335
$ws.Range("K84").Value = 72623.33499999999
$ws.Range("M84").Value = -67319.33499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 495.53333
$ws.Range("I107").Value = 444.41666
$ws.Range("K107").Value = 1333.24998
$ws.Range("M107").Value = 586.7500199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4399.8
$ws.Range("I132").Value = 4250
$ws.Range("K132").Value = 12750
$ws.Range("M132").Value = -10220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2960.4
$ws.Range("I136").Value = 2832
$ws.Range("J136").Value = 3233.25
$ws.Range("K136").Value = 8496
$ws.Range("L136").Value = 9699.75
$ws.Range("M136").Value = -5946
$ws.Range("N136").Value = -14799.75
